$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two discontinued products entirely (rows shift up).
# Delete the lower row first so the earlier row's number stays valid.
$ws.Rows.Item(123).EntireRow.Delete()
$ws.Rows.Item(114).EntireRow.Delete()
